$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values formatted as text (e.g. "1.007", "28.097.31").
# Excel auto-converts single-dot numeric-looking text to a Number on assignment,
# so force those specific cells to Text format first to preserve the literal string.

$ws.Range("D2").Value = '28.097.31'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '1.789.48'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.007'
$ws.Range("E4").Value = '  +0.47%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.68'
$ws.Range("E5").Value = '  -0.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5216'
$ws.Range("E7").Value = '  +2.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3798'
$ws.Range("E8").Value = '  -3.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07929'
$ws.Range("E9").Value = '  -3.52%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.42'
$ws.Range("E10").Value = '  -0.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.090'
$ws.Range("E11").Value = '  -1.32%  '
$ws.Range("B12").Value = 'BinanceUSD'
$ws.Range("C12").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.006'
$ws.Range("E12").Value = '  +0.42%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.236'
$ws.Range("E13").Value = '  -0.62%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.42'
$ws.Range("E14").Value = '  -2.53%  '
$ws.Range("D15").Value = '1.801.48'
$ws.Range("E15").Value = '  -0.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.272'
$ws.Range("E16").Value = '  -2.93%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.50'
$ws.Range("E17").Value = '  -0.93%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001084'
$ws.Range("E18").Value = '  -5.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06576'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.24'
$ws.Range("E21").Value = '  -2.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.950'
$ws.Range("E22").Value = '  -2.34%  '
$ws.Range("D23").Value = '28.185.89'
$ws.Range("E23").Value = '  -0.91%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.07'
$ws.Range("E24").Value = '  -1.99%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.262'
$ws.Range("E25").Value = '  -0.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.19'
$ws.Range("E26").Value = '  +3.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.42'
$ws.Range("E27").Value = '  -3.40%  '
$ws.Range("D28").Value = '2.004.06'
$ws.Range("E28").Value = '  -1.05%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.328'
$ws.Range("E29").Value = '  -2.91%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.76'
$ws.Range("E30").Value = '  -2.44%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1080'
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.049'
$ws.Range("E32").Value = '  -5.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.692'
$ws.Range("E33").Value = '  +1.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.502'
$ws.Range("E34").Value = '  -4.98%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07189'
$ws.Range("E35").Value = '  +2.28%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.20'
$ws.Range("E36").Value = '  +8.54%  '
$ws.Range("E37").Value = '  -1.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2136'
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.050'
$ws.Range("E39").Value = '  -2.98%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.568'
$ws.Range("E40").Value = '  -2.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6155'
$ws.Range("E41").Value = '  -1.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.162'
$ws.Range("E42").Value = '  -0.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.375'
$ws.Range("E43").Value = '  -2.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.28'
$ws.Range("E44").Value = '  -1.30%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.769'
$ws.Range("E45").Value = '  +0.78%  '
$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5931'
$ws.Range("E46").Value = '  +0.67%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '126.19'
$ws.Range("E47").Value = '  +0.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.211'
$ws.Range("E48").Value = '  +2.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.912'
$ws.Range("E49").Value = '  -3.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06782'
$ws.Range("E50").Value = '  -1.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '72.36'
$ws.Range("E51").Value = '  -2.24%  '
